$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Rv3025c", 5, "iscS Rv3025c", "FUNCTION: Catalyzes the removal of elemental sulfur from cysteine to produce alanine (Probable). Participates in the biosynthesis of metalloclusters by providing the inorganic sulfur required for Fe-S core formation. One acceptor is Whib3, on which this enzyme assembles a 4Fe-4S cluster. It can use both L-cysteine and L-selenocysteine as substrates. {ECO:0000269|PubMed:17609386, ECO:0000269|PubMed:24548275}.", 17)
    ,@(11, "Rv1956", 5, "higA1 higA Rv1956", "FUNCTION: Antitoxin component of an atypical, type II toxin-antitoxin chaperone (TAC) system. Upon expression in M.smegmatis neutralizes the effect of cognate toxin HigB1. Neutralization of HigB1 toxin in E.coli or M.marinum also requires SecB-like chaperone Rv1957, making this the first toxin-antitoxin chaperone (TAC) system. Antitoxin aggregation and degradation are prevented by the chaperone. {ECO:0000269|PubMed:21536872}.; FUNCTION: In M.tuberculosis represses expression of the Rv1954A-higB1-higA1-Rv1957 operon promoter but not that of the higB1-higA1-Rv1957 operon. {ECO:0000269|PubMed:20585061}.", 17)
    ,@(12, "Rv2454c", 5, "korB Rv2454c", "FUNCTION: Component of KG oxidoreductase (KOR) that catalyzes the CoA-dependent oxidative decarboxylation of 2-oxoglutarate (alpha-ketoglutarate, KG) to succinyl-CoA. Methyl viologen can act as electron acceptor in vitro; the physiologic electron acceptor is unknown. Is involved in the alternative TCA pathway that functions concurrently with fatty acid beta-oxidation. Since a growing body of evidence indicates that lipids (for example cholesterol and fatty acids) are a predominant growth substrate for M.tuberculosis during infection, flux through KOR likely represents an important step in intermediary metabolism in vivo. KOR-dependent decarboxylation of KG also appears to be an important source of CO(2) in M.tuberculosis metabolism. {ECO:0000269|PubMed:19936047}.", 17)
    ,@(13, "Rv0819", 5, "mshD Rv0819", "FUNCTION: Catalyzes the transfer of acetyl from acetyl-CoA to desacetylmycothiol (Cys-GlcN-Ins) to form mycothiol. {ECO:0000269|PubMed:12375100}.", 17)
    ,@(15, "Rv2455c", 4, "korA Rv2455c", "FUNCTION: Component of KG oxidoreductase (KOR) that catalyzes the CoA-dependent oxidative decarboxylation of 2-oxoglutarate (alpha-ketoglutarate, KG) to succinyl-CoA. Methyl viologen can act as electron acceptor in vitro; the physiologic electron acceptor is unknown. Is involved in the alternative TCA pathway that functions concurrently with fatty acid beta-oxidation. Since a growing body of evidence indicates that lipids (for example cholesterol and fatty acids) are a predominant growth substrate for M.tuberculosis during infection, flux through KOR likely represents an important step in intermediary metabolism in vivo. KOR-dependent decarboxylation of KG also appears to be an important source of CO(2) in M.tuberculosis metabolism. {ECO:0000269|PubMed:19936047}.", 17)
    ,@(16, "Rv2388c", 3, "hemW hemN Rv2388c MTCY253.33", "FUNCTION: Probably acts as a heme chaperone, transferring heme to an unknown acceptor. Binds one molecule of heme per monomer, possibly covalently (By similarity). Binds 1 [4Fe-4S] cluster. The cluster is coordinated with 3 cysteines and an exchangeable S-adenosyl-L-methionine (By similarity). {ECO:0000250|UniProtKB:P32131, ECO:0000250|UniProtKB:P52062}.", 17)
    ,@(17, "Rv1957", 3, "secBL Rv1957", "FUNCTION: Chaperone component of an atypical, type II toxin-antitoxin chaperone (TAC) system. Prevents antitoxin HigA1 aggregation in vitro at a 1:3 chaperone:antitoxin ratio, probably also protects antitoxin HigA1 from protease. Required for neutralization of toxin HigB1 upon ectopic expression in Mycobacterium marinum or E.coli. When expressed in E.coli complements a secB deletion, restores export of OmpA and MBP and inhibits aggregation of proOmpC although it is less efficient than endogenous SecB. Complements the general chaperone function of E.coli SecB less well. {ECO:0000269|PubMed:21536872}.", 17)
    ,@(18, "Rv2404c", 3, "lepA Rv2404c MTCY253.16", "FUNCTION: Required for accurate and efficient protein synthesis under certain stress conditions. May act as a fidelity factor of the translation reaction, by catalyzing a one-codon backward translocation of tRNAs on improperly translocated ribosomes. Back-translocation proceeds from a post-translocation (POST) complex to a pre-translocation (PRE) complex, thus giving elongation factor G a second chance to translocate the tRNAs correctly. Binds to ribosomes in a GTP-dependent manner. {ECO:0000255|HAMAP-Rule:MF_00071, ECO:0000269|PubMed:22578863}.", 17)
    ,@(19, "Rv1324", 2, "Rv1324 MTCY130.09", "", 17)
    ,@(20, "Rv3410c", 2, "guaB3 Rv3410c MTCY78.18", "FUNCTION: Has no inosine-5'-monophosphate dehydrogenase activity. {ECO:0000269|PubMed:21081761}.", 17)
    ,@(21, "Rv2847c", 2, "cysG Rv2847c", "", 17)
    ,@(22, "Rv0154c", 1, "fadE2 Rv0154c", "", 17)
    ,@(23, "Rv2690c", 1, "Rv2690c", "FUNCTION: Probable amino-acid or metabolite transport protein. {ECO:0000256|ARBA:ARBA00002249}.", 17)
    ,@(24, "Rv2393", 1, "che1 Rv2393", "", 17)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
